$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Add the new row of data (row 10)
$ws.Range("B10").Value = "Programmazione"
$ws.Range("C10").Value = "Unity"
$ws.Range("D10").Value = 0.0625
$ws.Range("E10").Value = "Unity con Mattia, spiegato debug"

# Update selection to match new active cell
$ws.Range("E10").Select()
